# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 132
$ws1.Range("F3").Value  = 330
$ws1.Range("F4").Value  = 422
$ws1.Range("F5").Value  = 1714
$ws1.Range("F11").Value = 4851
$ws1.Range("F17").Value = 173
$ws1.Range("F21").Value = 3782
$ws1.Range("F30").Value = 83
$ws1.Range("F31").Value = 572
$ws1.Range("F32").Value = 8
$ws1.Range("F33").Value = 26
$ws1.Range("F34").Value = 904
$ws1.Range("F35").Value = 2419
$ws1.Range("F36").Value = 424

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 132
$ws4.Range("F3").Value  = 330
$ws4.Range("F4").Value  = 422
$ws4.Range("F5").Value  = 1714
$ws4.Range("F11").Value = 4851
$ws4.Range("F17").Value = 173
$ws4.Range("F21").Value = 3782
$ws4.Range("F30").Value = 83
$ws4.Range("F31").Value = 572
$ws4.Range("F32").Value = 8
$ws4.Range("F34").Value = 26
$ws4.Range("F35").Value = 904
$ws4.Range("F36").Value = 2419
$ws4.Range("F37").Value = 424
